$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44308, 0, 1, 25.4323499491353),
    @(44309, 0, 1, 25.4323499491353),
    @(44310, 0, 1, 25.4323499491353),
    @(44311, 1, 1, 25.4323499491353),
    @(44312, 0, 1, 25.4323499491353)
)

$startRow = 234
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    # copy the format of the A column date cell from the row above (keeps
    # the same style index instead of minting a new one)
    $ws.Range("A" + ($row - 1)).Copy() | Out-Null
    $ws.Range("A" + $row).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$excel.CutCopyMode = $false
